# Add a new slide using the "Title and Content" layout (classic PpSlideLayout
# ppLayoutText = 2) right after the existing slide, then fill in its title
# and the four body bullet points.

$p = $ppt.ActivePresentation

$s = $p.Slides.Add(2, 2)

$title = $s.Shapes.Item(1).TextFrame.TextRange
$title.Text = "New slide title"

$body = $s.Shapes.Item(2).TextFrame.TextRange
$body.Text = "A"
[void]$body.InsertAfter("`rFew")
[void]$body.InsertAfter("`rBullet")
[void]$body.InsertAfter("`rPoints")
